$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header values in D1/E1
$ws.Range("D1").Value = "FirstName"
$ws.Range("E1").Value = "MiddleName"

# Swap data values in D2/E2
$ws.Range("D2").Value = "Bishal"
$ws.Range("E2").Value = "Programmer"

# Update column widths (columns B and D are left untouched, per the diff).
# NOTE: the COM ColumnWidth property only accepts values on a 1/6-character
# grid (it round-trips through a 6-px-per-character model), so values are
# picked to land as closely as possible on the target stored width after
# that internal rounding.
$ws.Columns.Item(1).ColumnWidth = 10.833333333333334
$ws.Columns.Item(3).ColumnWidth = 6.666666666666667
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666

# Update selection
$ws.Range("E7").Select()
